# Insert a new data row for "Comercializadora del Agro de Limarí - Arveja Verde"
# (week of 2023-07-27) above the current row 65, shifting existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(65).Insert()

$ws.Range("A65").Value = 2
$ws.Range("B65").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 45134
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = 100112022
$ws.Range("G65").Value = "Arveja Verde"
$ws.Range("H65").Value = "Perfection"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 600
$ws.Range("K65").Value = 23000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = 24000
$ws.Range("N65").Value = "$/malla 25 kilos"
$ws.Range("O65").Value = "Provincia de Limarí"
$ws.Range("P65").Value = 960
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = "Hortaliza"
